{"js": "// Replace the placeholder \"TBD\" office-hours text with the final schedule.\nconst body = context.document.body;\nconst results = body.search(\"TBD\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Tuesdays, Thursdays 9:50am-11:35am\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Replace the placeholder \"TBD\" office-hours text with the final schedule.\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"TBD\"\n$find.Replacement.Text = \"Tuesdays, Thursdays 9:50am-11:35am\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$find.Execute([ref]\"TBD\", [ref]$true, [ref]$true, $false, $false, $false, $true, $true, $false, \"Tuesdays, Thursdays 9:50am-11:35am\", 2)\n"}
